$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New shared "stats" query used for column C (Cases/Samples/Files tabs).
# ---------------------------------------------------------------------------
$statsQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.stage_of_disease IN ['Unknown']  
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# ---------------------------------------------------------------------------
# Cases-tab query (column B, row 2): same as before, plus a trailing
# `Cohort` column in the RETURN clause.
# ---------------------------------------------------------------------------
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE diag.stage_of_disease IN ['Unknown']
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

# ---------------------------------------------------------------------------
# Files-tab query (column B, row 4): same as before, minus the trailing
# `Study Code` column in the RETURN clause.
# ---------------------------------------------------------------------------
$filesQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.stage_of_disease IN ['Unknown']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis 
'@

# Write order below controls the order in which brand-new shared strings
# are appended to xl/sharedStrings.xml (new entries append in first-write
# order, right after the entries that are still in use). To reproduce the
# canonical ordering (stats query, then files query, then cases query)
# the stats query is written first, then the files query, then the cases
# query.

# Row 2/3/4, column C: shared "stats" query (first use -> new index right
# after the unchanged Samples query).
$ws.Range("C2").Value = $statsQuery
$ws.Range("C3").Value = $statsQuery
$ws.Range("C4").Value = $statsQuery

# Row 4 (FilesTab)
$ws.Range("B4").Value = $filesQuery

# Row 2 (CasesTab)
$ws.Range("B2").Value = $casesQuery

# ---------------------------------------------------------------------------
# Row heights shrank now that the long queries wrap across fewer lines.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

# Selection moved from B2 to C2.
$ws.Range("C2").Select() | Out-Null
